# DIV-4872 - added UserRole tab to the template
# To be able to show/hide tabs based on a user role
#
# Inserts a new "UserRole" column into the CaseTypeTab definition sheet
# (between the existing "TabFieldDisplayOrder" and "FieldShowCondition"
# columns) and makes CaseTypeTab the active/selected sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CaseTypeTab")

# Insert a new column at I, shifting the existing I:K columns (and their
# data/styles) one place to the right, becoming J:L.
$ws.Range("I1").EntireColumn.Insert()

# Populate the new column's header description and field name.
$ws.Range("I2").Value = "MaxLength: 100. No entry for role means no role restriction for that tab. Enter role on a single row per tab"
$ws.Range("I3").Value = "UserRole"

# Make CaseTypeTab the active sheet/tab and select cell I4, matching the
# updated selection/active-tab state captured in the workbook.
$ws.Activate()
$ws.Range("I4").Select()
